# Mise à jour de l'application
# Adds a new training-day column (DN) for date 46063 (2026-02-10) to the
# attendance sheet, filling in each player's attendance status and then
# copying the cell formatting from the previous day's column (DM). The
# existing summary formulas in columns B:J already span the whole
# K:VR/VS/VT/... ranges so they recalc automatically once the new data is
# written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header date for the new training day.
$ws.Cells.Item(1, 118).Value = 46063

# Per-player attendance for the new date (rows 16, 17 and 25 are left blank,
# and rows 12, 21, 23 get no new cell at all since those players' data
# already stopped earlier in the sheet).
$ws.Cells.Item(2, 118).Value = "P"   # Alban Rambaud
$ws.Cells.Item(3, 118).Value = "P"   # Jassim Assoul
$ws.Cells.Item(4, 118).Value = "P"   # Enzo Vita
$ws.Cells.Item(5, 118).Value = "P"   # Romain Thunet
$ws.Cells.Item(6, 118).Value = "B"   # Amine Taiar
$ws.Cells.Item(7, 118).Value = "P"   # Naim Ighbane
$ws.Cells.Item(8, 118).Value = "P"   # Hedi Nasri
$ws.Cells.Item(9, 118).Value = "P"   # Mattheo Haon
$ws.Cells.Item(10, 118).Value = "P"  # Maé Clavel
$ws.Cells.Item(11, 118).Value = "P"  # Levy Ndoutoume
$ws.Cells.Item(13, 118).Value = "P"  # Rayane Chayebi
$ws.Cells.Item(14, 118).Value = "P"  # Ilan Ihaddadene
$ws.Cells.Item(15, 118).Value = "P"  # Karahali Souaré
$ws.Cells.Item(18, 118).Value = "B"  # Emmanuel Valey
$ws.Cells.Item(19, 118).Value = "M"  # Jeremie Laurent
$ws.Cells.Item(20, 118).Value = "P"  # Sofiane Belle
$ws.Cells.Item(22, 118).Value = "P"  # Naim Dhib
$ws.Cells.Item(24, 118).Value = "B"  # Yoan Zouma
$ws.Cells.Item(26, 118).Value = "P"  # Omar Benyounes
$ws.Cells.Item(27, 118).Value = "P"  # Yoann Martelat
$ws.Cells.Item(28, 118).Value = "P"  # Malik Boussaid
$ws.Cells.Item(29, 118).Value = "P"  # Kamal Bafounta
$ws.Cells.Item(30, 118).Value = "P"  # Theo Owono
$ws.Cells.Item(31, 118).Value = "P"  # Mehdi Boussaid
$ws.Cells.Item(32, 118).Value = "P"  # Nathanael Beta

# Copy the previous column's formatting onto every newly-populated (or
# deliberately left blank) DN cell, row by row (a multi-area range copy
# doesn't paste reliably here).
$rows = @(1,2,3,4,5,6,7,8,9,10,11,13,14,15,16,17,18,19,20,22,24,25,26,27,28,29,30,31,32)
foreach ($r in $rows) {
    $ws.Range("DM" + $r).Copy()
    $ws.Range("DN" + $r).PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = $false

# Recalculate so the COUNTA/COUNTIF summary columns (B:J) pick up the
# newly-added data.
$excel.CalculateFullRebuild()

# Restore the view state recorded in the saved file: the previously active
# cell is re-selected.
$ws.Range("DP28").Select()
